# Applies the 2025-11-20 Betfair Back/Lay sheet update:
#   - two fixtures inserted before the existing Brazilian Serie A rows
#     (Danish 1st Division @14:30, Romanian Liga II @14:30)
#   - the three original Brazilian Serie A rows shift from 2-4 down to 4-6
#     and get refreshed odds
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 2; the old rows 2-4 slide down to 4-6
$ws.Rows(2).Insert() | Out-Null
$ws.Rows(2).Insert() | Out-Null

# Row-insert copies the header row's bold/bordered style into the new rows -
# clear that so the new data rows stay plain, like the rest of the data rows
$ws.Range("A2:AO3").ClearFormats()

# Force the Date column to text so "2025-11-20" is stored as a literal string
# rather than being auto-parsed into a date serial number
$ws.Range("B2:B3").NumberFormat = "@"

# --- Row 2: Danish 1st Division - Hobro vs Hvidovre ---
$ws.Range("A2").Value = "Danish 1st Division"
$ws.Range("B2").Value = "2025-11-20"
$ws.Range("C2").Value = "14:30:00"
$ws.Range("D2").Value = "Hobro"
$ws.Range("E2").Value = "Hvidovre"
$ws.Range("F2").Value = 3.05
$ws.Range("G2").Value = 3.8
$ws.Range("H2").Value = 2.22
$ws.Range("I2").Value = 2.48
$ws.Range("J2").Value = 3.5
$ws.Range("K2").Value = 4.2
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 3.85
$ws.Range("O2").Value = 1.26
$ws.Range("P2").Value = 2.14
$ws.Range("Q2").Value = 1.72
$ws.Range("R2").Value = 1.24
$ws.Range("S2").Value = 2.62
$ws.Range("T2").Value = 1.53
$ws.Range("U2").Value = 1.01
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 1.35
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 15
$ws.Range("Z2").Value = 20
$ws.Range("AA2").Value = 44
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 12.5
$ws.Range("AD2").Value = 16
$ws.Range("AE2").Value = 29
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 20
$ws.Range("AH2").Value = 23
$ws.Range("AI2").Value = 46
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 48
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# --- Row 3: Romanian Liga II - Concordia Chiajna vs Bihor Oradea ---
$ws.Range("A3").Value = "Romanian Liga II"
$ws.Range("B3").Value = "2025-11-20"
$ws.Range("C3").Value = "14:30:00"
$ws.Range("D3").Value = "Concordia Chiajna"
$ws.Range("E3").Value = "Bihor Oradea"
$ws.Range("F3").Value = 1.04
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 1.04
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 1.03
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.34
$ws.Range("O3").Value = 1.01
$ws.Range("P3").Value = 1.34
$ws.Range("Q3").Value = 1.02
$ws.Range("R3").Value = 1.16
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 1.04
$ws.Range("U3").Value = 1.04
$ws.Range("V3").Value = 1.01
$ws.Range("W3").Value = 1.01
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# --- Row 4 (was row 2): Brazilian Serie A - Juventude vs Cruzeiro MG - updated odds ---
$ws.Range("F4").Value = 4.7
$ws.Range("G4").Value = 5.3
$ws.Range("H4").Value = 1.83
$ws.Range("I4").Value = 1.87
$ws.Range("J4").Value = 3.65
$ws.Range("K4").Value = 3.95
$ws.Range("N4").Value = 3.3
$ws.Range("O4").Value = 1.37
$ws.Range("P4").Value = 1.78
$ws.Range("Q4").Value = 2.1
$ws.Range("S4").Value = 3.75
$ws.Range("U4").Value = 1.89
$ws.Range("V4").Value = 2.14
$ws.Range("W4").Value = 1.23
$ws.Range("X4").Value = 13.5
$ws.Range("Y4").Value = 8.2
$ws.Range("Z4").Value = 11
$ws.Range("AA4").Value = 980
$ws.Range("AB4").Value = 16
$ws.Range("AD4").Value = 12.5
$ws.Range("AG4").Value = 21
$ws.Range("AH4").Value = 23
$ws.Range("AI4").Value = 55
$ws.Range("AJ4").Value = 150
$ws.Range("AK4").Value = 75
$ws.Range("AL4").Value = 80
$ws.Range("AO4").Value = 15.5

# --- Row 5 (was row 3): Brazilian Serie A - Bahia vs Fortaleza EC - updated odds ---
$ws.Range("G5").Value = 1.56
$ws.Range("K5").Value = 5
$ws.Range("L5").Value = 1.27
$ws.Range("M5").Value = 1.04
$ws.Range("O5").Value = 1.23
$ws.Range("P5").Value = 2.26
$ws.Range("Q5").Value = 1.69
$ws.Range("R5").Value = 1.5
$ws.Range("S5").Value = 2.74
$ws.Range("T5").Value = 1.81
$ws.Range("W5").Value = 2.74
$ws.Range("X5").Value = 25
$ws.Range("Y5").Value = 29
$ws.Range("Z5").Value = 60
$ws.Range("AA5").Value = 210
$ws.Range("AB5").Value = 11.5
$ws.Range("AC5").Value = 12.5
$ws.Range("AD5").Value = 30
$ws.Range("AF5").Value = 12
$ws.Range("AG5").Value = 12
$ws.Range("AH5").Value = 27
$ws.Range("AI5").Value = 90
$ws.Range("AJ5").Value = 18
$ws.Range("AK5").Value = 19.5
$ws.Range("AL5").Value = 40
$ws.Range("AN5").Value = 8.8

# --- Row 6 (was row 4): Brazilian Serie A - Corinthians vs Sao Paulo - updated odds ---
$ws.Range("F6").Value = 2.44
$ws.Range("G6").Value = 2.58
$ws.Range("H6").Value = 3.45
$ws.Range("N6").Value = 2.7
$ws.Range("Q6").Value = 2.58
$ws.Range("S6").Value = 5.3
$ws.Range("W6").Value = 1.63

